# [Demo Project][Add Employee][Version 1][Add Employee Done, Write on Excel Done]
#
# Writes the newly-added employee records into the worksheet. Employee ID,
# First/Middle Name and Last Name are written starting at row 1 (replacing
# the old header labels) and continuing down one row per employee.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Employee ID is numeric-looking text ("4187", ...) and must be stored as a
# genuine text value (not auto-converted to a number). Route it through a
# text formula + paste-values-only so the stored cell keeps its default
# (unstyled) cell format, matching how the source app writes these cells.
function Set-TextValue($cell, [string]$val) {
    $escaped = $val.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$employees = @(
    @("4187", "AgR6vFirstName 71Fm2MiddleName", "AbZIhLastName"),
    @("5318", "QdquzFirstName wpmoAMiddleName", "3fDcELastName"),
    @("5012", "ZuHNAFirstName qMOFcMiddleName", "lgZgiLastName"),
    @("8784", "nnV3lFirstName AZz7dMiddleName", "JJITYLastName")
)

$row = 1
foreach ($employee in $employees) {
    Set-TextValue $ws.Cells.Item($row, 1) $employee[0]
    $ws.Cells.Item($row, 2).Value = $employee[1]
    $ws.Cells.Item($row, 3).Value = $employee[2]
    $row = $row + 1
}
